$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.079.78'
$ws.Range('E2').Value = '  +2.22%  '
$ws.Range('D3').Value = '1.654.36'
$ws.Range('E3').Value = '  +2.39%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '214.16'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.47%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.528'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.59%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '23.57'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +3.56%  '
$ws.Range('E9').Value = '  +2.20%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0616'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.66%  '
$ws.Range('E11').Value = '  -1.26%  '
$ws.Range('D12').Value = '1.889.60'
$ws.Range('E12').Value = '  +2.40%  '
$ws.Range('D13').Value = '1.657.56'
$ws.Range('E13').Value = '  +2.49%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.09'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.59%  '
$ws.Range('E15').Value = '  +3.60%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.92'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.59%  '
$ws.Range('D17').Value = '28.090.39'
$ws.Range('E17').Value = '  +2.16%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '233.65'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.13%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.71'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.63%  '
$ws.Range('E20').Value = '  +0.94%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.00'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.14%  '
$ws.Range('E22').Value = '  +5.48%  '
$ws.Range('E23').Value = '  +3.20%  '
$ws.Range('E24').Value = '  +3.89%  '
$ws.Range('E25').Value = '  +0.94%  '
$ws.Range('E26').Value = '  +1.44%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '15.83'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.05%  '
$ws.Range('E28').Value = '  +0.78%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.07%  '
$ws.Range('E30').Value = '  +1.54%  '
$ws.Range('E31').Value = '  +0.60%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.36'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.62%  '
$ws.Range('D33').Value = '1.455.86'
$ws.Range('E33').Value = '  -0.57%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.10'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.54%  '
$ws.Range('E35').Value = '  +2.61%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.33'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.39%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.894'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +4.21%  '
$ws.Range('E38').Value = '  +1.66%  '
$ws.Range('B39').Value = 'TrustWalletToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.929'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.40%  '
$ws.Range('B40').Value = 'ImmutableX'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.561'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.58%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '69.51'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.35%  '
$ws.Range('E42').Value = '  +3.33%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.00'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.07%  '
$ws.Range('B44').Value = 'mCoin'
$ws.Range('C44').Value = 'https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.46'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.17%  '
$ws.Range('B45').Value = 'RenderToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.83'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +6.55%  '
$ws.Range('B46').Value = 'MXToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.23'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.23%  '
$ws.Range('E47').Value = '  +3.19%  '
$ws.Range('D48').Value = '1.797.93'
$ws.Range('E48').Value = '  +2.20%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '89.17'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.98%  '
$ws.Range('E50').Value = '  +0.98%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0508'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.10%  '
